$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.285.22"
$ws.Range("E2").Value = "  +1.20%  "
$ws.Range("D3").Value = "1.832.37"
$ws.Range("E3").Value = "  +0.78%  "
$ws.Range("D4").Value = "'1.012"
$ws.Range("E4").Value = "  +0.95%  "
$ws.Range("D5").Value = "'314.60"
$ws.Range("E5").Value = "  +1.54%  "
$ws.Range("D6").Value = "'1.009"
$ws.Range("E6").Value = "  +0.74%  "
$ws.Range("E7").Value = "  +1.69%  "
$ws.Range("D8").Value = "'0.3687"
$ws.Range("E8").Value = "  +0.68%  "
$ws.Range("E9").Value = "  +0.97%  "
$ws.Range("D10").Value = "'0.8855"
$ws.Range("E10").Value = "  +1.56%  "
$ws.Range("D11").Value = "'20.48"
$ws.Range("E11").Value = "  +1.05%  "
$ws.Range("D12").Value = "1.899.32"
$ws.Range("E12").Value = "  +4.07%  "
$ws.Range("D13").Value = "'0.07327"
$ws.Range("E13").Value = "  +3.02%  "
$ws.Range("D14").Value = "'5.429"
$ws.Range("E14").Value = "  +0.72%  "
$ws.Range("D15").Value = "'94.05"
$ws.Range("E15").Value = "  +2.92%  "
$ws.Range("D16").Value = "'6.561"
$ws.Range("D17").Value = "'1.009"
$ws.Range("D18").Value = "'0.000008790"
$ws.Range("E18").Value = "  +1.09%  "
$ws.Range("D19").Value = "'1.010"
$ws.Range("E19").Value = "  +0.82%  "
$ws.Range("D20").Value = "27.549.45"
$ws.Range("D21").Value = "'14.77"
$ws.Range("E21").Value = "  +0.76%  "
$ws.Range("D22").Value = "'5.284"
$ws.Range("E22").Value = "  -0.20%  "
$ws.Range("E23").Value = "  +0.72%  "
$ws.Range("D24").Value = "2.098.68"
$ws.Range("E24").Value = "  +2.52%  "
$ws.Range("D25").Value = "'1.889"
$ws.Range("E25").Value = "  -0.28%  "
$ws.Range("D26").Value = "'151.87"
$ws.Range("E26").Value = "  +0.54%  "
$ws.Range("D27").Value = "'18.64"
$ws.Range("E27").Value = "  +1.03%  "
$ws.Range("D28").Value = "'2.143"
$ws.Range("E28").Value = "  +0.11%  "
$ws.Range("D29").Value = "'5.227"
$ws.Range("E29").Value = "  -0.75%  "
$ws.Range("D30").Value = "'116.99"
$ws.Range("E30").Value = "  +0.35%  "
$ws.Range("D31").Value = "'0.08989"
$ws.Range("E31").Value = "  +1.02%  "
$ws.Range("D32").Value = "'0.7497"
$ws.Range("E32").Value = "  -1.41%  "
$ws.Range("D33").Value = "'1.174"
$ws.Range("E33").Value = "  +0.73%  "
$ws.Range("D34").Value = "'4.543"
$ws.Range("E34").Value = "  +1.18%  "
$ws.Range("D35").Value = "'2.940"
$ws.Range("E35").Value = "  +1.29%  "
$ws.Range("D36").Value = "'1.010"
$ws.Range("E36").Value = "  +0.86%  "
$ws.Range("D37").Value = "'1.089"
$ws.Range("E37").Value = "  -0.18%  "
$ws.Range("D38").Value = "'0.05343"
$ws.Range("E38").Value = "  +0.98%  "
$ws.Range("D39").Value = "'0.01956"
$ws.Range("E39").Value = "  +0.41%  "
$ws.Range("D40").Value = "'2.974"
$ws.Range("E40").Value = "  +0.12%  "
$ws.Range("D41").Value = "'2.391"
$ws.Range("E41").Value = "  +2.60%  "
$ws.Range("D42").Value = "'7.242"
$ws.Range("E42").Value = "  +1.01%  "
$ws.Range("D43").Value = "'0.5289"
$ws.Range("E43").Value = "  +0.01%  "
$ws.Range("E44").Value = "  +0.08%  "
$ws.Range("D45").Value = "'8.483"
$ws.Range("D46").Value = "'0.4930"
$ws.Range("E46").Value = "  +1.44%  "
$ws.Range("D47").Value = "'10.49"
$ws.Range("E47").Value = "  +0.10%  "
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").Value = "'104.95"
$ws.Range("E48").Value = "  +1.46%  "
$ws.Range("B49").Value = "PaxDollar"
$ws.Range("C49").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D49").Value = "'1.010"
$ws.Range("E49").Value = "  +0.78%  "
$ws.Range("D50").Value = "'1.667"
$ws.Range("E50").Value = "  +0.30%  "
$ws.Range("E51").Value = "  -0.01%  "
